$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.225.58"
$ws.Range("E2").Value = "  -3.05%  "

$ws.Range("D3").Value = "3.135.06"
$ws.Range("E3").Value = "  -2.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.45%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "3.135.18"
$ws.Range("E8").Value = "  -2.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.153"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.478"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.09%  "

$ws.Range("D15").Value = "3.644.49"
$ws.Range("E15").Value = "  -2.25%  "

$ws.Range("D16").Value = "64.260.08"
$ws.Range("E16").Value = "  -3.19%  "

$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").Value = "3.131.95"
$ws.Range("E18").Value = "  -2.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.39%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  -2.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.46%  "

$ws.Range("E29").Value = "  -3.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.33%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.50%  "

$ws.Range("E35").Value = "  -5.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.67%  "

$ws.Range("D39").Value = "0.0₃0754"
$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "450.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.11%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.124"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.09%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0401"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.68%  "

$ws.Range("D44").Value = "2.882.24"
$ws.Range("E44").Value = "  -1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.273"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("E49").Value = "  -1.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.02%  "

$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
